$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last refreshed" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 11 de Julio de 2020 a las 13:54"

# Row 5 - Brasil: new totals (country ranking unaffected)
$ws.Range("B5").Value = 1807496
$ws.Range("C5").Value = 3158
$ws.Range("E5").Value = 523383
$ws.Range("G5").Value = 77
$ws.Range("H5").Value = 70601

# Row 24 - Catar
$ws.Range("B24").Value = 103128
$ws.Range("C24").Value = 498
$ws.Range("D24").Value = 98934
$ws.Range("E24").Value = 4048

# Row 32 - Bielorrusia
$ws.Range("B32").Value = 64767
$ws.Range("C32").Value = 163
$ws.Range("D32").Value = 54919
$ws.Range("E32").Value = 9389
$ws.Range("G32").Value = 5
$ws.Range("H32").Value = 459

# Row 48 - Afganistan
$ws.Range("B48").Value = 34366
$ws.Range("C48").Value = 172
$ws.Range("D48").Value = 21135
$ws.Range("E48").Value = 12237
$ws.Range("G48").Value = 23
$ws.Range("H48").Value = 994

# Row 64 - Nepal
$ws.Range("B64").Value = 16719
$ws.Range("C64").Value = 70
$ws.Range("D64").Value = 8442
$ws.Range("E64").Value = 8239
$ws.Range("G64").Value = 3
$ws.Range("H64").Value = 38

# Row 70 - Uzbekistan
$ws.Range("D70").Value = 7540
$ws.Range("E70").Value = 4610
$ws.Range("G70").Value = 2
$ws.Range("H70").Value = 56

# Rows 74/75 - Kenia overtakes Australia in ranking; country names swap,
# row 74 gets fresh Kenia data, row 75 gets Australia's prior (row 74) data
$ws.Range("A74").Value = "Kenia"
$ws.Range("B74").Value = 9726
$ws.Range("C74").Value = 278
$ws.Range("D74").Value = 2832
$ws.Range("E74").Value = 6710
$ws.Range("G74").Value = 3
$ws.Range("H74").Value = 184

$ws.Range("A75").Value = "Australia"
$ws.Range("B75").Value = 9549
$ws.Range("C75").Value = 190
$ws.Range("D75").Value = 7730
$ws.Range("E75").Value = 1712
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 107

# Rows 97/98/99 - Madagascar overtakes Republica de Africa Central and
# Hungria; labels shift down, row 97 gets fresh Madagascar data, rows
# 98/99 inherit the prior rows' (97/98) data
$ws.Range("A97").Value = "Madagascar"
$ws.Range("B97").Value = 4578
$ws.Range("C97").Value = 435
$ws.Range("D97").Value = 2287
$ws.Range("E97").Value = 2257
$ws.Range("H97").Value = 34

$ws.Range("A98").Value = "Republica de Africa Central"
$ws.Range("B98").Value = 4259
$ws.Range("C98").Value = 0
$ws.Range("D98").Value = 1142
$ws.Range("E98").Value = 3064
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 53

$ws.Range("A99").Value = "Hungria"
$ws.Range("B99").Value = 4229
$ws.Range("C99").Value = 6
$ws.Range("D99").Value = 2974
$ws.Range("E99").Value = 660
$ws.Range("G99").Value = 2
$ws.Range("H99").Value = 595

# Row 120 - Islandia
$ws.Range("B120").Value = 1888
$ws.Range("C120").Value = 2
$ws.Range("D120").Value = 1860
$ws.Range("E120").Value = 18

# Row 138 - Burkina Faso
$ws.Range("D138").Value = 864
$ws.Range("E138").Value = 103

Write-Output "applied edits"
